# Merge the three runs of the first paragraph in "TextBox 82" (shape id 83,
# the slide's only top-level shape named "TextBox 82") into a single run,
# without altering the visible text.
#
# The engine's TextRange.Text setter diffs old vs. new content and keeps
# any runs whose text is unchanged, so assigning the exact same text back
# is a no-op for the run layout. To force a real merge we first overwrite
# the paragraph with placeholder text that shares no common wording with
# the original (so none of the old runs can be matched/preserved), which
# collapses the paragraph down to one run; we then set the real text,
# which again lands in a single run. The placeholder is chosen to be
# close in length to the final text so the auto-fit text box does not
# need to reflow/resize in between (avoiding any incidental size change).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$sh = $s.Shapes.Item("TextBox 82")

$finalText = "To the left are two distinctive sonic signatures that one might compare using the tool " + [char]8211 + " Rosalind and Claudio. At a glance you can see quickly that one has dark bars above the line and one has dark bars below the line."

$tr = $sh.TextFrame.TextRange
$para1 = $tr.Paragraphs(1, 1)

$placeholder = "bicycle zeppelin lantern crimson crimson horizon bicycle obelisk marmalade phantom glacier quixotic zeppelin marmalade thunder crimson cascade gravel zeppelin obelisk thunder obelisk glacier crimson tundra phantom lantern"
$para1.Text = $placeholder

$para1b = $tr.Paragraphs(1, 1)
$para1b.Text = $finalText
